# Update crypto price/volume data per commit: "Updated cryptos list on Mon Oct 30 16:53:08 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.571.93'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.816.64'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.50'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.579'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.99%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '34.99'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +7.27%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.301'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.50%  '
$ws.Range('E10').Value = '  +0.73%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '2.078.95'
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = '1.805.20'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.647'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').Value = '34.530.47'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.47'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '245.82'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.58'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +3.80%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('E23').Value = '  +0.72%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '171.99'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.56%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.10'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  +4.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.81'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.92%  '
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +2.47%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.25'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0530'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.76%  '
$ws.Range('E33').Value = '  +1.07%  '
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D35').Value = '1.403.57'
$ws.Range('E35').Value = '  -1.09%  '
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('E37').Value = '  +1.75%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '83.27'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.86'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +4.28%  '
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.41'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  +1.58%  '
$ws.Range('E45').Value = '  +2.80%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0512'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.96%  '
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '1.978.99'
$ws.Range('E48').Value = '  +1.67%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '105.56'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0130'
$ws.Range('E51').Value = '  -0.77%  '
